$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new log entries by Bruno Díaz (2017-05-27) for the task
# "Sprint 3 - Integración BackEnd y FrontEnd", added below the existing
# rows (row 51 was the last one, now rows 52 and 53 are appended).

$xlPasteFormats = -4122

$ws.Range("A52").Value = "Bruno Díaz"
$ws.Range("B52").Value = 42882
$ws.Range("C52").Value = 3
$ws.Range("D52").Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Range("E52").Value = "Se crearon 5 ejemplos de llamadas a servicios (1 DELETE, 1 POST, 1 PUT, 2 GET -Uno con y uno sin parámetro-)"

$ws.Range("A53").Value = "Bruno Díaz"
$ws.Range("B53").Value = 42882
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Range("E53").Value = "Se modificaron lo 5 ejemplos para pasar a trabajar con objetos dentro de objetos."

# Reuse the existing date format/style from the row above (B51) instead of
# creating a brand new number format entry.
$ws.Range("B51").Copy()
$ws.Range("B52:B53").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("E54").Select()
